# Weekly update: a new week of price data for "Acelga" at Vega Monumental
# Concepción was inserted into the daily log. The new record is inserted
# right after the existing row 45 pair (same chronological slot the source
# system uses), pushing every subsequent record down by one pair of rows
# (Primera / Segunda quality rows), and the sheet grows from A1:R157 to
# A1:R159.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new pair of rows (Primera + Segunda) by inserting two
# blank rows at 46:47; everything at/after row 46 shifts down by two rows.
$ws.Rows("46:47").Insert()

# Seed the new rows with the same template values as the adjacent existing
# record (same market/region/category/quality layout), then set the new
# date for this week's entry.
$ws.Range("A44:R45").Copy($ws.Range("A46"))

$ws.Cells.Item(46, 4).Value = 44469
$ws.Cells.Item(47, 4).Value = 44469
